# Update document date and the 25 division-problem answers to the
# new revision (commit: "Update master to output generated at c8c62b6").

$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-15 Friday", 2) | Out-Null

$d.Content.Find.Execute("617÷6=102, 5", $true, $false, $false, $false, $false, $true, 1, $false, "651÷6=108, 3", 2) | Out-Null
$d.Content.Find.Execute("642÷9=71, 3", $true, $false, $false, $false, $false, $true, 1, $false, "453÷7=64, 5", 2) | Out-Null
$d.Content.Find.Execute("669÷8=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "434÷3=144, 2", 2) | Out-Null
$d.Content.Find.Execute("439÷5=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "486÷7=69, 3", 2) | Out-Null
$d.Content.Find.Execute("777÷3=259, 0", $true, $false, $false, $false, $false, $true, 1, $false, "989÷6=164, 5", 2) | Out-Null

$d.Content.Find.Execute("485÷7=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "441÷4=110, 1", 2) | Out-Null
$d.Content.Find.Execute("389÷3=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "578÷7=82, 4", 2) | Out-Null
$d.Content.Find.Execute("540÷8=67, 4", $true, $false, $false, $false, $false, $true, 1, $false, "818÷3=272, 2", 2) | Out-Null
$d.Content.Find.Execute("265÷6=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "682÷2=341, 0", 2) | Out-Null
$d.Content.Find.Execute("778÷5=155, 3", $true, $false, $false, $false, $false, $true, 1, $false, "891÷5=178, 1", 2) | Out-Null

$d.Content.Find.Execute("722÷9=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "358÷5=71, 3", 2) | Out-Null
$d.Content.Find.Execute("689÷9=76, 5", $true, $false, $false, $false, $false, $true, 1, $false, "595÷4=148, 3", 2) | Out-Null
$d.Content.Find.Execute("621÷9=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "839÷8=104, 7", 2) | Out-Null
$d.Content.Find.Execute("793÷9=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "124÷5=24, 4", 2) | Out-Null
$d.Content.Find.Execute("281÷4=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "710÷5=142, 0", 2) | Out-Null

$d.Content.Find.Execute("838÷3=279, 1", $true, $false, $false, $false, $false, $true, 1, $false, "956÷5=191, 1", 2) | Out-Null
$d.Content.Find.Execute("971÷7=138, 5", $true, $false, $false, $false, $false, $true, 1, $false, "858÷3=286, 0", 2) | Out-Null
$d.Content.Find.Execute("247÷5=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "114÷7=16, 2", 2) | Out-Null
$d.Content.Find.Execute("419÷5=83, 4", $true, $false, $false, $false, $false, $true, 1, $false, "470÷2=235, 0", 2) | Out-Null
$d.Content.Find.Execute("842÷9=93, 5", $true, $false, $false, $false, $false, $true, 1, $false, "574÷4=143, 2", 2) | Out-Null

$d.Content.Find.Execute("766÷3=255, 1", $true, $false, $false, $false, $false, $true, 1, $false, "721÷3=240, 1", 2) | Out-Null
$d.Content.Find.Execute("497÷2=248, 1", $true, $false, $false, $false, $false, $true, 1, $false, "236÷3=78, 2", 2) | Out-Null
$d.Content.Find.Execute("465÷8=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "292÷6=48, 4", 2) | Out-Null
$d.Content.Find.Execute("702÷7=100, 2", $true, $false, $false, $false, $false, $true, 1, $false, "354÷4=88, 2", 2) | Out-Null
$d.Content.Find.Execute("902÷2=451, 0", $true, $false, $false, $false, $false, $true, 1, $false, "524÷5=104, 4", 2) | Out-Null

Write-Output "Replacements complete."
